$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), carrying formatting/styles along
$ws.Columns("D").Insert()

# Copy formatting from the (now shifted) former column D, which is now column E, into new column D
# (done per contiguous data block so rows that never had a D:K cell, e.g. section header rows,
# do not get a stray empty D cell created)
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate new column D with the latest quarter figures
$ws.Range("D7").Value = 43407
$ws.Range("D8").Value = 95400
$ws.Range("D9").Value = 61700
$ws.Range("D10").Value = 33700
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 14400
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 118400
$ws.Range("D18").Value = -23000
$ws.Range("D20").Value = 100
$ws.Range("D21").Value = -16300
$ws.Range("D22").Value = 100
$ws.Range("D23").Value = -23000
$ws.Range("D24").Value = -6700
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = -16200
$ws.Range("D27").Value = -16200
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -100
$ws.Range("D33").Value = -16200
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = -16200
$ws.Range("D38").Value = 43407
$ws.Range("D41").Value = 10700
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 17100
$ws.Range("D44").Value = 40400
$ws.Range("D45").Value = 10900
$ws.Range("D46").Value = 79100
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 79800
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 20500
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 179500
$ws.Range("D57").Value = 37400
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 12300
$ws.Range("D60").Value = 49700
$ws.Range("D61").Value = 0
$ws.Range("D62").Value = 35000
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 84700
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 141500
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 94800
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43407
$ws.Range("D81").Value = -16200
$ws.Range("D83").Value = 6600
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = -5100
$ws.Range("D91").Value = -7400
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -7400
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -12600
